# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status cells move from "In Translation" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime / Name get the new handback's values
#  - The stale "version mismatch" error detail is cleared
#  - The now-longer status/name columns are widened to fit their new content

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("L2").Value = "2016-12-09 06:33:28"
$zhcn.Range("M2").Value = "TestHandback_201612090233"
$zhcn.Range("R2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(13).ColumnWidth = 27.166666666666668
$zhcn.Columns.Item(18).ColumnWidth = 12.833333333333334

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("L2").Value = "2016-12-09 06:33:46"
$dede.Range("M2").Value = "TestHandback_201612090233"
$dede.Range("R2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(13).ColumnWidth = 27.166666666666668
$dede.Columns.Item(18).ColumnWidth = 12.833333333333334
